$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 2
$ws.Range("F6").Value = -6
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 2
